# Added Min and Max Values for CO2 Plant Size
# The restrictions-table columns for Usage (pkm / tkm-SZM / tkm-N1 / tkm-N2 / tkm-N3)
# were reordered. Re-label the headers in AE1:AI1 and move each column's
# data value in row 2 so it keeps following its (renamed) header.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Re-order the Usage headers (columns AE..AI on row 1) ---
# Old layout: AE=pkmUsage       AF=tkm-SZMUsage  AG=tkm-N1Usage  AH=tkm-N2Usage  AI=tkm-N3Usage
# New layout: AE=tkm-N2Usage    AF=tkm-N3Usage   AG=tkm-N1Usage  AH=pkmUsage     AI=tkm-SZMUsage
# (AG1/AG2 = tkm-N1Usage / 7.5 stay exactly where they were)
$ws.Range("AE1").Value = "tkm-N2Usage"
$ws.Range("AF1").Value = "tkm-N3Usage"
$ws.Range("AH1").Value = "pkmUsage"
$ws.Range("AI1").Value = "tkm-SZMUsage"

# --- Move the row-2 values so each one still sits under its own header ---
$ws.Range("AE2").Value = 24.2
$ws.Range("AF2").Value = 130.3
$ws.Range("AH2").Value = 858
$ws.Range("AI2").Value = 414.5

# --- Tiny floating point recalculation drift on a handful of other cells ---
$ws.Range("E2").Value = 612.0528792959244
$ws.Range("H2").Value = 741.7852278688525
$ws.Range("M2").Value = 514.1244186085765
$ws.Range("N2").Value = 514.1244186085765
$ws.Range("R2").Value = 288.4396604831752
